# Updated cryptos list on Wed Nov  1 07:21:24 UTC 2023 with GitHub Actions
#
# Note: several "Price" values in column D look numeric (e.g. "224.85")
# but must stay plain text, matching the source data (which also keeps
# clearly non-numeric values like "34.479.25" as text). Assigning such a
# numeric-looking string straight to .Value would make Excel silently
# convert it to a real number (and mangle precision). Prefixing with a
# leading apostrophe forces Excel to keep it as literal text, exactly like
# typing '224.85 into a cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "34.479.25"
$ws.Range("E2").Value = "  +0.38%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.807.50"
$ws.Range("E3").Value = "  +0.36%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.20%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'224.85"
$ws.Range("E5").Value = "  -1.04%  "

# Row 6 - XRP
$ws.Range("D6").Value = "'0.594"
$ws.Range("E6").Value = "  +3.37%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.17%  "

# Row 8 - Solana
$ws.Range("D8").Value = "'38.32"
$ws.Range("E8").Value = "  +6.69%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -3.78%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.0672"
$ws.Range("E10").Value = "  -2.90%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.0977"
$ws.Range("E11").Value = "  +1.26%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "2.067.55"
$ws.Range("E12").Value = "  +0.30%  "

# Row 13 - Chainlink
$ws.Range("D13").Value = "'11.09"
$ws.Range("E13").Value = "  -4.65%  "

# Row 14 - WrappedEther
$ws.Range("D14").Value = "1.808.15"
$ws.Range("E14").Value = "  +0.05%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "'0.629"
$ws.Range("E15").Value = "  -2.09%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "34.446.31"
$ws.Range("E16").Value = "  +0.27%  "

# Row 17 - Polkadot
$ws.Range("D17").Value = "'4.37"

# Row 18 - Litecoin
$ws.Range("D18").Value = "'68.02"
$ws.Range("E18").Value = "  -1.42%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "'241.71"
$ws.Range("E19").Value = "  -1.29%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0769"
$ws.Range("E20").Value = "  -3.14%  "

# Row 21 - Avalanche
$ws.Range("D21").Value = "'11.09"
$ws.Range("E21").Value = "  -4.21%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  -0.12%  "

# Row 23 - Uniswap
$ws.Range("D23").Value = "'4.10"
$ws.Range("E23").Value = "  -1.41%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +1.13%  "

# Row 25 - Monero
$ws.Range("D25").Value = "'170.99"
$ws.Range("E25").Value = "  -0.68%  "

# Row 26 - Cosmos
$ws.Range("D26").Value = "'7.70"
$ws.Range("E26").Value = "  -3.15%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "'17.37"
$ws.Range("E27").Value = "  +3.11%  "

# Row 28 - Stellar
$ws.Range("D28").Value = "'0.120"
$ws.Range("E28").Value = "  +1.52%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  -0.21%  "

# Row 30 - PancakeSwap
$ws.Range("D30").Value = "'1.23"
$ws.Range("E30").Value = "  -1.36%  "

# Row 31 - Filecoin
$ws.Range("E31").Value = "  -1.88%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  -4.06%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "'0.0512"
$ws.Range("E33").Value = "  -3.23%  "

# Row 34 - LidoDAOToken
$ws.Range("D34").Value = "'1.82"
$ws.Range("E34").Value = "  +0.18%  "

# Row 35 - Maker
$ws.Range("D35").Value = "1.317.09"
$ws.Range("E35").Value = "  -5.75%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  -4.53%  "

# Row 37 - TrustWalletToken
$ws.Range("D37").Value = "'1.06"
$ws.Range("E37").Value = "  -0.97%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  -0.88%  "

# Row 39 - RenderToken
$ws.Range("E39").Value = "  -6.16%  "

# Row 40 - Aave
$ws.Range("D40").Value = "'82.99"
$ws.Range("E40").Value = "  +0.45%  "

# Row 41 - HuobiToken
$ws.Range("E41").Value = "  +0.82%  "

# Row 42 and 43 swap: WEMIXToken <-> MXToken (plus new values)
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "'2.81"
$ws.Range("E42").Value = "  -0.46%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'1.21"
$ws.Range("E43").Value = "  -1.08%  "

# Row 44 - ARBITRUM
$ws.Range("D44").Value = "'0.944"
$ws.Range("E44").Value = "  -1.79%  "

# Row 45 - InjectiveProtocol
$ws.Range("E45").Value = "  +1.75%  "

# Row 46 - Kaspa
$ws.Range("D46").Value = "'0.0512"
$ws.Range("E46").Value = "  +1.09%  "

# Row 47 - RocketPoolETH
$ws.Range("D47").Value = "1.968.26"
$ws.Range("E47").Value = "  +0.31%  "

# Row 48 - FraxShare
$ws.Range("E48").Value = "  -4.58%  "

# Row 49 - PaxDollar
$ws.Range("E49").Value = "  -0.16%  "

# Row 50 - Quant
$ws.Range("D50").Value = "'102.28"
$ws.Range("E50").Value = "  -1.84%  "

# Row 51 - BabyDogeCoin
$ws.Range("E51").Value = "  -7.15%  "
